$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15 (shifts existing rows 15-23 down to 16-24)
$ws.Rows("15:15").Insert()

# Populate the new row with the CLC population entry
$ws.Range("A15").Value = "CLC"
$ws.Range("B15").Value = "Snake River Coho Salmon"
$ws.Range("C15").Value = "Clearwater River"
$ws.Range("D15").Value = "CRSFC-c"
$ws.Range("E15").Value = "South Fork Clearwater River"

Write-Output ($ws.Range("A14").Text + "|" + $ws.Range("A15").Text + "|" + $ws.Range("A16").Text + "|" + $ws.Range("A24").Text)
